# Fixed issues with 81RF protective element:
# Changed default xls parameters to disable 81x protections
# (81RFRP, 81RFDFP, 81RF Trip Delay columns T:V and AH:AJ, rows 2-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

for ($r = 2; $r -le 15; $r++) {
    $ws.Range("T$r").Value = 100
    $ws.Range("U$r").Value = 10
    $ws.Range("V$r").Value = 0.1

    $ws.Range("AH$r").Value = 100
    $ws.Range("AI$r").Value = 10
    $ws.Range("AJ$r").Value = 0.1
}

# Update view/selection to mirror the author's workbook state
$ws.Activate() | Out-Null
$ws.Range("AH3:AJ15").Select() | Out-Null
